$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "Subsi" + hidden _GoBack bookmark + "dies" -> single "Subsidies"
# run. Replacing the (contiguous) found text collapses the two runs into one
# and consumes the bookmark that used to sit between them.
# ---------------------------------------------------------------------------
$subsidiesRange = $d.Range(0, $d.Content.End)
$subsidiesRange.Find.Execute("Subsidies", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Subsidies", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 2: drop the paragraph that only holds a manual page break, right
# before the "Results" subtitle, and move the hidden _GoBack bookmark so it
# marks the (now merged) start of the "Results" paragraph.
# ---------------------------------------------------------------------------
function Find-ResultsParagraphIndex($doc) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text.TrimEnd() -eq "Results" -and $p.Style.NameLocal -eq "Subtitle") {
            return $i
        }
    }
    return -1
}

$resultsIdx = Find-ResultsParagraphIndex $d
if ($resultsIdx -gt 1) {
    $pageBreakPara = $d.Paragraphs.Item($resultsIdx - 1)
    $pageBreakPara.Range.Delete() | Out-Null
}

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete() | Out-Null
}

$resultsIdx2 = Find-ResultsParagraphIndex $d
if ($resultsIdx2 -gt 0) {
    $resultsPara = $d.Paragraphs.Item($resultsIdx2)
    $goBackRange = $d.Range($resultsPara.Range.Start, $resultsPara.Range.Start)
    $d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null
}
